# Digitale_Kompetenzen.xlsx - "Fixed Assessment Type" edit
#
# The AssessmentType column (BC) currently stores a plain numeric flag (1)
# for every data row. Replace it with the proper assessment-type label
# "SOFTSKILL" and give those cells the small accent font used for type
# labels elsewhere in the workbook (8pt "MesloLGM NF", custom gold color,
# vertically centered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2..56 (row 1 is the header row).
$firstRow = 2
$lastRow = 56
$col = "BC"

$rng = $ws.Range("$col$firstRow`:$col$lastRow")

# New assessment-type label (was a bare numeric 1 in every row).
$rng.Value = "SOFTSKILL"

# Apply the accent font/format used for the corrected AssessmentType cells.
$rng.Font.Name = "MesloLGM NF"
$rng.Font.Size = 8
$rng.Font.Color = 7185097   # BGR for RGB FFC9A26D
$rng.VerticalAlignment = -4108   # xlCenter

# Keep the selection/visible area close to where the edit happened.
$ws.Activate()
$ws.Range("BD60").Select()
